$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 137,2
$data[0,0] = "good child"
$data[0,1] = "いい子|いいこ"
$data[1,0] = "color"
$data[1,1] = "色|いろ"
$data[2,0] = "boxed lunch"
$data[2,1] = "お弁当|おべんとう"
$data[3,0] = "Kabuki; traditional Japanese theatrical art"
$data[3,1] = "歌舞伎|かぶき"
$data[4,0] = "guitar"
$data[4,1] = "ギター"
$data[5,0] = "last year"
$data[5,1] = "去年|きょねん"
$data[6,0] = "medicine"
$data[6,1] = "薬|くすり"
$data[7,0] = "to take medicine"
$data[7,1] = "薬を飲む|くすりをのむ"
$data[8,0] = "concert"
$data[8,1] = "コンサート"
$data[9,0] = "near future"
$data[9,1] = "今度|こんど"
$data[10,0] = "essay; composition"
$data[10,1] = "作文|さくぶん"
$data[11,0] = "exam"
$data[11,1] = "試験|しけん"
$data[12,0] = "ski"
$data[12,1] = "スキー"
$data[13,0] = "last month"
$data[13,1] = "先月|せんげつ"
$data[14,0] = "word; vocabulary"
$data[14,1] = "単語|たんご"
$data[15,0] = "piano"
$data[15,1] = "ピアノ"
$data[16,0] = "pizza"
$data[16,1] = "ピザ"
$data[17,0] = "illness; sickness"
$data[17,1] = "病気|びょうき"
$data[18,0] = "blue"
$data[18,1] = "青い|あおい"
$data[19,0] = "red"
$data[19,1] = "赤い|あかい"
$data[20,0] = "black"
$data[20,1] = "黒い|くろい"
$data[21,0] = "lonely"
$data[21,1] = "寂しい|さびしい"
$data[22,0] = "white"
$data[22,1] = "白い|しろい"
$data[23,0] = "young"
$data[23,1] = "若い|わかい"
$data[24,0] = "mean-spirited"
$data[24,1] = "意地悪|いじわる(な)"
$data[25,0] = "to dance"
$data[25,1] = "踊る|おどる"
$data[26,0] = "(something) ends"
$data[26,1] = "終わる|おわる"
$data[27,0] = "to be popular"
$data[27,1] = "人気がある|にんきがある"
$data[28,0] = "(something) begins"
$data[28,1] = "始まる|はじまる"
$data[29,0] = "to play (a string instrument or piano)"
$data[29,1] = "弾く|ひく"
$data[30,0] = "to get (from somebody)"
$data[30,1] = "もらう"
$data[31,0] = "to memorize"
$data[31,1] = "覚える|おぼえる"
$data[32,0] = "to appear; to attend; to exit"
$data[32,1] = "出る|でる"
$data[33,0] = "to do physical exercises"
$data[33,1] = "運動する|うんどうする"
$data[34,0] = "to take a walk"
$data[34,1] = "散歩する|さんぽする"
$data[35,0] = "from..."
$data[35,1] = "～から"
$data[36,0] = "by all means"
$data[36,1] = "是非|ぜひ"
$data[37,0] = "by the way"
$data[37,1] = "ところで"
$data[38,0] = "all"
$data[38,1] = "みんな"
$data[39,0] = "already"
$data[39,1] = "もう"
$data[40,0] = "one"
$data[40,1] = "一つ|ひとつ"
$data[41,0] = "two"
$data[41,1] = "二つ|ふたつ"
$data[42,0] = "three"
$data[42,1] = "三つ|みっつ"
$data[43,0] = "four"
$data[43,1] = "四つ|よっつ"
$data[44,0] = "five"
$data[44,1] = "五つ|いつつ"
$data[45,0] = "six"
$data[45,1] = "六つ|むっつ"
$data[46,0] = "seven"
$data[46,1] = "七つ|ななつ"
$data[47,0] = "eight"
$data[47,1] = "八つ|やっつ"
$data[48,0] = "nine"
$data[48,1] = "九つ|ここのつ"
$data[49,0] = "ten"
$data[49,1] = "十|とお"
$data[50,0] = "black"
$data[50,1] = "黒い|くろい"
$data[51,0] = "white"
$data[51,1] = "白い|しろい"
$data[52,0] = "red"
$data[52,1] = "赤い|あかい"
$data[53,0] = "blue"
$data[53,1] = "青い|あおい"
$data[54,0] = "yellow"
$data[54,1] = "黄色い|きいろい"
$data[55,0] = "brown"
$data[55,1] = "茶色い|ちゃいろい"
$data[56,0] = "green"
$data[56,1] = "緑／グリーン|みどり／---"
$data[57,0] = "purple"
$data[57,1] = "紫|むらさき"
$data[58,0] = "gray"
$data[58,1] = "灰色／グレー|はいいろ／---"
$data[59,0] = "light blue"
$data[59,1] = "水色|みずいろ"
$data[60,0] = "pink"
$data[60,1] = "ピンク"
$data[61,0] = "gold"
$data[61,1] = "金色／ゴールド|きんいろ／---"
$data[62,0] = "silver"
$data[62,1] = "銀色／シルバー|ぎんいろ／---"
$data[63,0] = "color"
$data[63,1] = "色|いろ"
$data[64,0] = "red bag"
$data[64,1] = "赤いかばん"
$data[65,0] = "I like red the best."
$data[65,1] = "赤がいちばん好きです。"
$data[66,0] = "black cat"
$data[66,1] = "黒猫|くろねこ"
$data[67,0] = "green sweater"
$data[67,1] = "緑／グリーンのセーター"
$data[68,0] = "Your face looks pale."
$data[68,1] = "顔が青いですね。|かおがあおいですね。"
$data[69,0] = "black and white picture"
$data[69,1] = "白黒の写真|しろくろのしゃしん"
$data[70,0] = "Mary has blonde hair."
$data[70,1] = "メアリーさんは金髪です。|メアリーさんはきんぱつです。"
$data[71,0] = "A.M."
$data[71,1] = "午前|ごぜん"
$data[72,0] = "P.M.; in the afternoon"
$data[72,1] = "午後|ごご"
$data[73,0] = "in the morning"
$data[73,1] = "午前中|ごぜんちゅう"
$data[74,0] = "noon"
$data[74,1] = "正午|しょうご"
$data[75,0] = "after..."
$data[75,1] = "～の後|～のあと"
$data[76,0] = "later"
$data[76,1] = "後で|あとで"
$data[77,0] = "back; behind"
$data[77,1] = "後ろ|うしろ"
$data[78,0] = "lastly"
$data[78,1] = "最後に|さいごに"
$data[79,0] = "before; front"
$data[79,1] = "前|まえ"
$data[80,0] = "name"
$data[80,1] = "名前|なまえ"
$data[81,0] = "advance sale"
$data[81,1] = "前売り|まえうり"
$data[82,0] = "to speak"
$data[82,1] = "話す|はなす"
$data[83,0] = "talk; story"
$data[83,1] = "話|はなし"
$data[84,0] = "telephone"
$data[84,1] = "電話|でんわ"
$data[85,0] = "conversation"
$data[85,1] = "会話|かいわ"
$data[86,0] = "little"
$data[86,1] = "少し|すこし"
$data[87,0] = "few"
$data[87,1] = "少ない|すくない"
$data[88,0] = "a little"
$data[88,1] = "少々|しょうしょう"
$data[89,0] = "girl"
$data[89,1] = "少女|しょうじょ"
$data[90,0] = "boy"
$data[90,1] = "少年|しょうねん"
$data[91,0] = "A.M."
$data[91,1] = "午前|ごぜん"
$data[92,0] = "P.M.; in the afternoon"
$data[92,1] = "午後|ごご"
$data[93,0] = "in the morning"
$data[93,1] = "午前中|ごぜんちゅう"
$data[94,0] = "noon"
$data[94,1] = "正午|しょうご"
$data[95,0] = "after..."
$data[95,1] = "～の後|～のあと"
$data[96,0] = "later"
$data[96,1] = "後で|あとで"
$data[97,0] = "back; behind"
$data[97,1] = "後ろ|うしろ"
$data[98,0] = "lastly"
$data[98,1] = "最後に|さいごに"
$data[99,0] = "before; front"
$data[99,1] = "前|まえ"
$data[100,0] = "name"
$data[100,1] = "名前|なまえ"
$data[101,0] = "advance sale"
$data[101,1] = "前売り|まえうり"
$data[102,0] = "time"
$data[102,1] = "時間|じかん"
$data[103,0] = "two hours"
$data[103,1] = "二時間|にじかん"
$data[104,0] = "between"
$data[104,1] = "間|あいだ"
$data[105,0] = "human being"
$data[105,1] = "人間|にんげん"
$data[106,0] = "one week"
$data[106,1] = "一週間|いっしゅうかん"
$data[107,0] = "house"
$data[107,1] = "家|いえ"
$data[108,0] = "family"
$data[108,1] = "家族|かぞく"
$data[109,0] = "house; home"
$data[109,1] = "家|うち"
$data[110,0] = "my wife"
$data[110,1] = "家内|かない"
$data[111,0] = "author"
$data[111,1] = "作家|さっか"
$data[112,0] = "to speak"
$data[112,1] = "話す|はなす"
$data[113,0] = "talk; story"
$data[113,1] = "話|はなし"
$data[114,0] = "telephone"
$data[114,1] = "電話|でんわ"
$data[115,0] = "conversation"
$data[115,1] = "会話|かいわ"
$data[116,0] = "little"
$data[116,1] = "少し|すこし"
$data[117,0] = "few"
$data[117,1] = "少ない|すくない"
$data[118,0] = "a little"
$data[118,1] = "少々|しょうしょう"
$data[119,0] = "girl"
$data[119,1] = "少女|しょうじょ"
$data[120,0] = "boy"
$data[120,1] = "少年|しょうねん"
$data[121,0] = "old (for things)"
$data[121,1] = "古い|ふるい"
$data[122,0] = "secondhand"
$data[122,1] = "中古|ちゅうこ"
$data[123,0] = "ancient times"
$data[123,1] = "古代|こだい"
$data[124,0] = "to know"
$data[124,1] = "知る|しる"
$data[125,0] = "acquaintance (formal)"
$data[125,1] = "知人|ちじん"
$data[126,0] = "acquaintance"
$data[126,1] = "知り合い|しりあい"
$data[127,0] = "to come"
$data[127,1] = "来る|くる"
$data[128,0] = "to come (long-form)"
$data[128,1] = "来ます|きます"
$data[129,0] = "not to come"
$data[129,1] = "来ない|こない"
$data[130,0] = "next week"
$data[130,1] = "来週|らいしゅう"
$data[131,0] = "visit to Japan"
$data[131,1] = "来日|らいにち"
$data[132,0] = "Thank you for everything."
$data[132,1] = "いろいろおせわになりました。"
$data[133,0] = "Please take care of yourself."
$data[133,1] = "体に気をつけてください。"
$data[134,0] = "I am looking forward to seeing you."
$data[134,1] = "お会いできるのを楽しみにしています。"
$data[135,0] = "Congratulations on..."
$data[135,1] = "～おめでとう（ございます）。"
$data[136,0] = "Happy birthday."
$data[136,1] = "（お）たんじょうびおめでとう。"

$ws.Range("A2:B138").Value = $data
